$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date placeholder text from
#    2/17/2018 -> 4/5/2019 everywhere it appears: the slide master, every
#    slide layout, and the notes master.
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "2/17/2018") {
                $tr.Text = "4/5/2019"
            }
        }
    }
}

# Slide master.
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout hanging off the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShapes $layouts.Item($li).Shapes
}

# Notes master.
Update-DateShapes $p.NotesMaster.Shapes

# ---------------------------------------------------------------------------
# 2) Fix the typo "peelk" -> "peek" in slide 4's speaker notes.
# ---------------------------------------------------------------------------
$notesShapes = $p.Slides.Item(4).NotesPage.Shapes
for ($i = 1; $i -le $notesShapes.Count; $i++) {
    $sh = $notesShapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -like "*peelk*") {
            $tr.Text = $tr.Text.Replace("peelk", "peek")
        }
    }
}
